$wb = $excel.ActiveWorkbook

# ----- Sheet "Overall": update row 2 stats -----
$overall = $wb.Worksheets.Item("Overall")
$overall.Range("B2").Value = 46
$overall.Range("C2").Value = 17
$overall.Range("D2").Value = 1.3368588490326596
$overall.Range("E2").Value = 0.41851851851851862
$overall.Range("F2").Value = 1.4307800191988784
$overall.Range("G2").Value = 29
$overall.Range("H2").Value = 18
$overall.Range("I2").Value = 47
$overall.Range("J2").Value = 348
$overall.Range("K2").Value = 23

# ----- Sheet "Zones": update per-zone stats -----
$zones = $wb.Worksheets.Item("Zones")

# Row 2 (Zone 1)
$zones.Range("B2").Value = 2
$zones.Range("C2").Value = 1
$zones.Range("D2").Value = 1.1111111111111112
$zones.Range("F2").Value = 1.1111111111111112

# Row 3 (Zone 2)
$zones.Range("B3").Value = 4
$zones.Range("D3").Value = 0.86041666666666661
$zones.Range("F3").Value = 0.86041666666666661

# Row 4 (Zone 3)
$zones.Range("B4").Value = 6
$zones.Range("D4").Value = 1.5333333333333334
$zones.Range("E4").Value = 0.88333333333333341
$zones.Range("F4").Value = 1.7933333333333334

# Row 5 (Zone 4)
$zones.Range("B5").Value = 3
$zones.Range("D5").Value = 1.8250000000000002
$zones.Range("F5").Value = 1.8250000000000002

# Row 6 (Zone 5)
$zones.Range("B6").Value = 1
$zones.Range("C6").Value = 4
$zones.Range("D6").Value = 2.6791666666666667
$zones.Range("F6").Value = 2.6791666666666667

# Row 7 (Zone 6)
$zones.Range("B7").Value = 5
$zones.Range("D7").Value = 0.81388888888888899
$zones.Range("E7").Value = 0.21666666666666667
$zones.Range("F7").Value = 0.93333333333333335

# Row 8 (Zone 7)
$zones.Range("B8").Value = 2
$zones.Range("C8").Value = 1
$zones.Range("D8").Value = 0.4333333333333334
$zones.Range("F8").Value = 0.4333333333333334

# Row 9 (Zone 8)
$zones.Range("B9").Value = 4
$zones.Range("C9").Value = 4
$zones.Range("D9").Value = 1.6181818181818179
$zones.Range("E9").Value = 0.1333333333333333
$zones.Range("F9").Value = 1.7666666666666664

# Row 10 (Zone 9)
$zones.Range("B10").Value = 4
$zones.Range("D10").Value = 0.72500000000000009
$zones.Range("E10").Value = 0.56666666666666687
$zones.Range("F10").Value = 0.77777777777777779

# Row 11 (Zone 10)
$zones.Range("C11").Value = 1
$zones.Range("D11").Value = 1.6062500000000002
$zones.Range("F11").Value = 1.6062500000000002

# Row 12 (Zone 11)
$zones.Range("B12").Value = 5
$zones.Range("D12").Value = 0.71333333333333337
$zones.Range("E12").Value = 0.33333333333333348
$zones.Range("F12").Value = 0.96666666666666667

# Row 13 (Zone 12)
$zones.Range("B13").Value = 6
$zones.Range("C13").Value = 6
$zones.Range("D13").Value = 1.392156862745098
$zones.Range("F13").Value = 1.392156862745098

# Row 14 (Zone 13)
$zones.Range("B14").Value = 4
$zones.Range("D14").Value = 0.82380952380952377
$zones.Range("E14").Value = 0.20833333333333337
$zones.Range("F14").Value = 1.0699999999999998
